$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

$ws.Cells.Item(21,3).Value = "138_卡罗拉_Carola_Rosa rugosa Thunb._20stems"
$ws.Cells.Item(21,6).Value = "'15"
$ws.Cells.Item(22,3).Value = "148_坦尼克_Tineke_Rosa rugosa Thunb._20stems"
$ws.Cells.Item(22,6).Value = "'15"
$ws.Cells.Item(23,3).Value = "148_坦尼克_Tineke_Rosa rugosa Thunb._20stems"
$ws.Cells.Item(23,6).Value = "'12"
$ws.Cells.Item(24,3).Value = "624_多丁白_undefined_undefined_1bunch"
$ws.Cells.Item(24,6).Value = "'10"
$ws.Cells.Item(25,3).Value = "626_多丁黄言_undefined_undefined_undefinedundefined"
$ws.Cells.Item(25,6).Value = "'10"
$ws.Cells.Item(26,1).Value = "'4"
$ws.Cells.Item(26,3).Value = "574_迷你菊白_undefined_undefined_1bunch"
$ws.Cells.Item(26,6).Value = "'15"
$ws.Cells.Item(27,3).Value = "576_迷你菊紫_undefined_undefined_1bunch"
$ws.Cells.Item(27,6).Value = "'5"
$ws.Cells.Item(28,3).Value = "575_迷你菊深粉_undefined_undefined_1bunch"
$ws.Cells.Item(28,6).Value = "'5"
$ws.Cells.Item(29,3).Value = "277_草莓杏仁饼_undefined_Rosa rugosa Thunb._10stems"
$ws.Cells.Item(29,6).Value = "'10"
$ws.Cells.Item(30,3).Value = "276_情迷罗拉_undefined_Rosa rugosa Thunb._10stems"
$ws.Cells.Item(30,6).Value = "'18"
$ws.Cells.Item(31,3).Value = "274_仙子之吻_undefined_Rosa rugosa Thunb._10stems"
$ws.Cells.Item(31,6).Value = "'9"
$ws.Cells.Item(32,3).Value = "225_果汁阳台_Juicy Terrazza_Rosa rugosa Thunb._10stems"
$ws.Cells.Item(32,6).Value = "'15"
$ws.Cells.Item(33,3).Value = "625_多丁紫蝴蝶_undefined_undefined_1bunch"
$ws.Cells.Item(33,6).Value = "'10"
$ws.Cells.Item(34,1).Value = "'5"
$ws.Cells.Item(34,3).Value = "144_高原红_High Plateau Red_Rosa rugosa Thunb._20stems"
$ws.Cells.Item(34,6).Value = "'20"
$ws.Cells.Item(35,3).Value = "147_娜欧米_Red Naomi_Rosa rugosa Thunb._20stems"
$ws.Cells.Item(35,6).Value = "'5"
$ws.Cells.Item(36,3).Value = "170_奶油杯_Butter Cup_Rosa rugosa Thunb._20stems"
$ws.Cells.Item(36,6).Value = "'7"
$ws.Cells.Item(37,3).Value = "479_绿灵草_lepidium_undefined_1bunch"
$ws.Cells.Item(37,6).Value = "'12"
$ws.Cells.Item(38,1).Value = "'"
$ws.Cells.Item(38,3).Value = "154_莫泊_Moab_Rosa rugosa Thunb._20stems"
$ws.Cells.Item(38,6).Value = "'8"
$ws.Cells.Item(39,1).Value = "'6"
$ws.Cells.Item(39,3).Value = "147_娜欧米_Red Naomi_Rosa rugosa Thunb._20stems"
$ws.Cells.Item(39,6).Value = "'8"
$ws.Cells.Item(40,3).Value = "181_月光女神_undefined_Rosa rugosa Thunb._20stems"
$ws.Cells.Item(40,6).Value = "'11"
$ws.Cells.Item(41,3).Value = "160_卡布奇诺_Cappuccino_Rosa rugosa Thunb._20stems"
$ws.Cells.Item(41,6).Value = "'9"
$ws.Cells.Item(42,3).Value = "135_甜蜜曼塔_sweet menta_Rosa rugosa Thunb._20stems"
$ws.Cells.Item(42,6).Value = "'10"
$ws.Cells.Item(43,3).Value = "203_佛罗伊德_Floyd_Rosa rugosa Thunb._20stems"
$ws.Cells.Item(43,6).Value = "'4"
$ws.Cells.Item(44,3).Value = "479_绿灵草_lepidium_undefined_1bunch"
$ws.Cells.Item(44,6).Value = "'13"
$ws.Cells.Item(45,1).Value = "'7"
$ws.Cells.Item(45,3).Value = "203_佛罗伊德_Floyd_Rosa rugosa Thunb._20stems"
$ws.Cells.Item(45,6).Value = "'5"
$ws.Cells.Item(46,3).Value = "170_奶油杯_Butter Cup_Rosa rugosa Thunb._20stems"
$ws.Cells.Item(46,6).Value = "'1"
$ws.Cells.Item(47,3).Value = "412_紫罗兰粉_violet pink_undefined_1bunch"
$ws.Cells.Item(47,6).Value = "'15"
$ws.Cells.Item(48,3).Value = "5_绿洋桔梗_Light Green Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
$ws.Cells.Item(48,6).Value = "'5"
$ws.Cells.Item(49,3).Value = "12_肉粉洋桔梗_Peach Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
$ws.Cells.Item(49,6).Value = "'30"
$ws.Cells.Item(50,3).Value = "1_白洋桔梗_White Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
$ws.Cells.Item(50,6).Value = "'20"
$ws.Cells.Item(51,3).Value = "509_翠珠粉_Didiscus caeruleus`npink_Trachymene Coerulea_1bunch"

$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(2,7).Value = "'0201020105555530105151030103015151515121010155510189151020571288119104135115530200"
